$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.283.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.171.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.72%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.170.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.01%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.727.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.349.18"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "419.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.01"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.487"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.71"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.49%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.31"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.15"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691.69"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.85%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.35"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.06"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.709"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.64%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0620"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.61"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "292.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.98%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0988"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.68%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -14.67%  "
